# Update countries & provincias Spain
# - Refresh COVID case numbers for several countries.
# - Re-rank rows whose "Casos totales" (col B) moved past a neighboring
#   country: Armenia now outranks Austria, and Hungria now outranks both
#   Ghana and Estado de Palestina. The non-updated countries keep their own
#   data, they just shift down a row as the updated country leapfrogs them.
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - refreshed totals -----------------------------
$ws.Range("B4").Value = 8388012
$ws.Range("C4").Value = 213
$ws.Range("D4").Value = 5457912
$ws.Range("E4").Value = 2705368
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 224732

# --- Rows 60-61: Armenia leapfrogs Austria ---------------------------------
# Row 60 becomes Armenia with freshly updated figures.
$ws.Range("A60").Value = "Armenia"
$ws.Range("B60").Value = 65460
$ws.Range("C60").Value = 766
$ws.Range("D60").Value = 48208
$ws.Range("E60").Value = 16161
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 1091

# Row 61 becomes Austria, carrying its previous (unchanged) figures.
$ws.Range("A61").Value = "Austria"
$ws.Range("B61").Value = 64806
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 49561
$ws.Range("E61").Value = 14352
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 893

# --- Rows 71-73: Hungria leapfrogs Ghana and Estado de Palestina ----------
# Row 71 becomes Hungria with freshly updated figures.
$ws.Range("A71").Value = "Hungria"
$ws.Range("B71").Value = 47768
$ws.Range("C71").Value = 1478
$ws.Range("D71").Value = 14312
$ws.Range("E71").Value = 32283
$ws.Range("G71").Value = 31
$ws.Range("H71").Value = 1173

# Row 72 becomes Ghana, carrying its previous (unchanged) figures.
$ws.Range("A72").Value = "Ghana"
$ws.Range("B72").Value = 47310
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 46618
$ws.Range("E72").Value = 382
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 310

# Row 73 becomes Estado de Palestina, carrying its previous (unchanged) figures.
$ws.Range("A73").Value = "Estado de Palestina"
$ws.Range("B73").Value = 47135
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 40498
$ws.Range("E73").Value = 6229
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 408

# --- Row 83: refreshed totals ----------------------------------------------
$ws.Range("E83").Value = 3740
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 926

# --- Row 178: refreshed totals ---------------------------------------------
$ws.Range("B178").Value = 540
$ws.Range("C178").Value = 5
$ws.Range("D178").Value = 493
$ws.Range("E178").Value = 40

# --- Timestamp footer (cell A1, above the table) ---------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 09:16"
